# Update "Generate Report for Handback" timestamps across sheets.
# These cells hold the date/time as text (shared strings) with a
# custom display format already applied via cell style, so we set
# the value explicitly as text to avoid Excel reinterpreting it as
# a date serial number.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for f5d6eaa3 row
$overview.Range("G4").Value = "2016-08-13 02:55:32"

# zh-cn sheet: Correspond Handoff / Handback datetimes for f5d6eaa3 row
$zhcn.Range("H4").Value = "2016-08-13 02:55:25"
$zhcn.Range("K4").Value = "2016-08-13 02:55:54"

# de-de sheet: Correspond Handoff / Handback datetimes for f5d6eaa3 row
$dede.Range("H4").Value = "2016-08-13 02:55:32"
$dede.Range("K4").Value = "2016-08-13 02:56:09"
